# ---------------------------------------------------------------------------
# "outreg dist. data 18-20" - add 2018/2019/2020 inflow+otflow columns
# (new G:L), pushing the existing 2021/2022 data (old C:F) right to I:L.
# Re-uses the existing header/data cell styling (centered; bold for the
# inflow/otflow + year header rows and for the district-name column).
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Row 2: inflow/otflow header labels, now spanning C:L ---
foreach ($col in @(3,5,7,9,11)) {
    $ws.Cells.Item(2, $col).Value = "inflow"
}
foreach ($col in @(4,6,8,10,12)) {
    $ws.Cells.Item(2, $col).Value = "otflow"
}
$hdrRow2 = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item(2, 12))
$hdrRow2.HorizontalAlignment = $xlCenter
$hdrRow2.Font.Bold = $true

# --- Row 3: year labels, now 2018,2018,2019,2019,2020,2020,2021,2021,2022,2022 ---
$years = @(2018,2018,2019,2019,2020,2020,2021,2021,2022,2022)
for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Cells.Item(3, 3 + $i).Value = $years[$i]
}
$hdrRow3 = $ws.Range($ws.Cells.Item(3, 3), $ws.Cells.Item(3, 12))
$hdrRow3.HorizontalAlignment = $xlCenter
$hdrRow3.Font.Bold = $true

# --- Rows 4-12: district data, columns C..L; row label (col B) becomes bold ---
$row4 = @(124177,80436,103076,68599,90215,62860,95323,65345,94031,64066)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4[$i]
}
$ws.Range($ws.Cells.Item(4, 3), $ws.Cells.Item(4, 12)).HorizontalAlignment = $xlCenter
$lbl4 = $ws.Cells.Item(4, 2)
$lbl4.HorizontalAlignment = $xlCenter
$lbl4.Font.Bold = $true

$row5 = @(9608,5978,8413,5352,7302,4995,7094,4806,6929,5047)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $row5[$i]
}
$ws.Range($ws.Cells.Item(5, 3), $ws.Cells.Item(5, 12)).HorizontalAlignment = $xlCenter
$lbl5 = $ws.Cells.Item(5, 2)
$lbl5.HorizontalAlignment = $xlCenter
$lbl5.Font.Bold = $true

$row6 = @(69556,51953,64646,50775,57076,45421,63328,48569,62611,47291)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, 3 + $i).Value = $row6[$i]
}
$ws.Range($ws.Cells.Item(6, 3), $ws.Cells.Item(6, 12)).HorizontalAlignment = $xlCenter
$lbl6 = $ws.Cells.Item(6, 2)
$lbl6.HorizontalAlignment = $xlCenter
$lbl6.Font.Bold = $true

$row7 = @(5143,2795,4991,2748,4325,2870,4174,2784,4362,2649)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, 3 + $i).Value = $row7[$i]
}
$ws.Range($ws.Cells.Item(7, 3), $ws.Cells.Item(7, 12)).HorizontalAlignment = $xlCenter
$lbl7 = $ws.Cells.Item(7, 2)
$lbl7.HorizontalAlignment = $xlCenter
$lbl7.Font.Bold = $true

$row8 = @(3182,1473,3065,1453,2649,1343,2328,1291,2313,1432)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, 3 + $i).Value = $row8[$i]
}
$ws.Range($ws.Cells.Item(8, 3), $ws.Cells.Item(8, 12)).HorizontalAlignment = $xlCenter
$lbl8 = $ws.Cells.Item(8, 2)
$lbl8.HorizontalAlignment = $xlCenter
$lbl8.Font.Bold = $true

$row9 = @(9618,3706,9107,3784,7546,3686,7081,3582,7017,3359)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, 3 + $i).Value = $row9[$i]
}
$ws.Range($ws.Cells.Item(9, 3), $ws.Cells.Item(9, 12)).HorizontalAlignment = $xlCenter
$lbl9 = $ws.Cells.Item(9, 2)
$lbl9.HorizontalAlignment = $xlCenter
$lbl9.Font.Bold = $true

$row10 = @(3669,1291,3663,1301,3122,1376,3038,1321,2932,1313)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, 3 + $i).Value = $row10[$i]
}
$ws.Range($ws.Cells.Item(10, 3), $ws.Cells.Item(10, 12)).HorizontalAlignment = $xlCenter
$lbl10 = $ws.Cells.Item(10, 2)
$lbl10.HorizontalAlignment = $xlCenter
$lbl10.Font.Bold = $true

$row11 = @(5303,1827,5272,1816,4724,1864,4650,1689,4440,1703)
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, 3 + $i).Value = $row11[$i]
}
$ws.Range($ws.Cells.Item(11, 3), $ws.Cells.Item(11, 12)).HorizontalAlignment = $xlCenter
$lbl11 = $ws.Cells.Item(11, 2)
$lbl11.HorizontalAlignment = $xlCenter
$lbl11.Font.Bold = $true

$row12 = @(4175,1258,3919,1370,3471,1305,3630,1303,3427,1272)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, 3 + $i).Value = $row12[$i]
}
$ws.Range($ws.Cells.Item(12, 3), $ws.Cells.Item(12, 12)).HorizontalAlignment = $xlCenter
$lbl12 = $ws.Cells.Item(12, 2)
$lbl12.HorizontalAlignment = $xlCenter
$lbl12.Font.Bold = $true

# --- Selection, as recorded after the edit ---
$null = $ws.Range("D18").Select()
